$d = $word.ActiveDocument
$apos = [char]39
$rsq = [char]8217

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the document to right after
#    the run "...resource can't be booked" (end of that paragraph).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("resource can" + $rsq + "t be booked", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(0)
$r1.InsertAfter([char]1)
$null = $d.Bookmarks.Add("_GoBack", $r1)
$r1.Text = ""

# ---------------------------------------------------------------------------
# 2. "Dispatcher drags a work order from work order list into a technicia's
#    schedule in the gantt chart" -> "...into a  technician's schedule in the
#    gantt chart"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("Dispatcher drags a work order from work order list into a technicia" + $apos + "s schedule in the gantt chart", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "Dispatcher drags a work order from work order list into a  technician" + $apos + "s schedule in the gantt chart"

# ---------------------------------------------------------------------------
# 3. "System validates changes to the gantt chart matching resource
#    availability and work order requirements" -> "...matching  technician's
#    availability and work order requirements"
# ---------------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("System validates changes to the gantt chart matching resource availability and work order requirements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Text = "System validates changes to the gantt chart matching  technician" + $apos + "s availability and work order requirements"

# ---------------------------------------------------------------------------
# 4. "Dispatcher clicks cancel on a work order in a resource's schedule"
#    (2nd occurrence only) -> "...in a technician's schedule"
# ---------------------------------------------------------------------------
$needle4 = "Dispatcher clicks cancel on a work order in a resource" + $apos + "s schedule"
$r4a = $d.Content
$null = $r4a.Find.Execute($needle4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4b = $d.Range($r4a.End, $d.Content.End)
$null = $r4b.Find.Execute($needle4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4b.Text = "Dispatcher clicks cancel on a work order in a technician" + $apos + "s schedule"

# ---------------------------------------------------------------------------
# 5. " System displays Work Order still in resource's schedule"
#    (2nd occurrence only) -> "... still in  technician's schedule"
# ---------------------------------------------------------------------------
$needle5 = "System displays Work Order still in resource" + $apos + "s schedule"
$r5a = $d.Content
$null = $r5a.Find.Execute($needle5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5b = $d.Range($r5a.End, $d.Content.End)
$null = $r5b.Find.Execute($needle5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5b.Text = "System displays Work Order still in  technician" + $apos + "s schedule"

# ---------------------------------------------------------------------------
# 6. "2.1 System unschedules resource from work order" (2nd occurrence only)
#    -> "2.1 System un-schedules  technician's from work order"
# ---------------------------------------------------------------------------
$needle6 = "2.1 System unschedules resource from work order"
$r6a = $d.Content
$null = $r6a.Find.Execute($needle6, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r6b = $d.Range($r6a.End, $d.Content.End)
$null = $r6b.Find.Execute($needle6, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r6b.Text = "2.1 System un-schedules  technician" + $apos + "s from work order"

# ---------------------------------------------------------------------------
# 7. "2.3 System displays schedule with work order no longer scheduled to
#    resource" -> "...scheduled to  technician"
# ---------------------------------------------------------------------------
$r7 = $d.Content
$null = $r7.Find.Execute("2.3 System displays schedule with work order no longer scheduled to resource", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r7.Text = "2.3 System displays schedule with work order no longer scheduled to  technician"

Write-Host "done"
